# classifications.xlsx - "Update wrt 2.3.0.dev1-steel" (Harpprecht et al. 2025)
# Adds ~47 new ecoinvent-style classification rows for steel / pig iron / direct
# reduction / electrowinning processes, matching the author's two-step entry:
# first rows 3109-3142 (after which the AutoFilter was (re)applied), then the
# remaining rows 3143-3155 appended without re-applying the filter - which is
# exactly why the final AutoFilter / _FilterDatabase range stops at row 3142
# even though the sheet's used range grows to row 3155.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRowsPhase1 = @(
    ,@("steel production, blast furnace-basic oxygen furnace, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, blast furnace-basic oxygen furnace, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("alloys production, for low-alloyed steel","alloys, for low-alloyed steel","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("pig iron production, blast furnace, with carbon capture and storage","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("carbon dioxide, captured at pig iron production plant, using monoethanolamine","carbon dioxide, captured","2011:Manufacture of basic chemicals","34210: Hydrogen, nitrogen, oxygen, carbon dioxide and rare gases; inorganic oxygen compounds of non-metals n.e.c.")
    ,@("steel production, blast furnace-basic oxygen furnace, low-alloyed, with carbon capture and storage","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, blast furnace-basic oxygen furnace, unalloyed, with carbon capture and storage","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("pig iron production, top gas recycling-blast furnace","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, low-alloyed, with top gas recycling","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, blast furnace-basic oxygen furnace, unalloyed, with top gas recycling","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("pig iron production, blast furnace, with top gas recycling, with carbon capture and storage","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("carbon dioxide, captured at steel production plant, using vacuum pressure swing adsorption","carbon dioxide, captured","2011:Manufacture of basic chemicals","34210: Hydrogen, nitrogen, oxygen, carbon dioxide and rare gases; inorganic oxygen compounds of non-metals n.e.c.")
    ,@("steel production, blast furnace-basic oxygen furnace, with top gas recycling, low-alloyed, with carbon capture and storage","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, blast furnace-basic oxygen furnace, with top gas recycling, unalloyed, with carbon capture and storage","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("iron porduction, with natural gas-based direct reduction","iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, natural gas-based direct reduction iron-electric arc furnace, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, natural gas-based direct reduction iron-electric arc furnace, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("iron production, with natural gas-based direct reduction, with carbon capture and storage","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("carbon dioxide, captured at steel production plant using direct reduction iron, using vacuum pressure swing adsorption","carbon dioxide, captured","2011:Manufacture of basic chemicals","34210: Hydrogen, nitrogen, oxygen, carbon dioxide and rare gases; inorganic oxygen compounds of non-metals n.e.c.")
    ,@("steel production, natural gas-based direct reduction iron-electric arc furnace, low-alloyed, with carbon capture and storage","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, natural gas-based direct reduction iron-electric arc furnace, unalloyed, with carbon capture and storage","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, hydrogen-based direct reduction iron-electric arc furnace, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, hydrogen-based direct reduction iron-electric arc furnace, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("pig iron production, hydrogen-based direct reduction iron","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("preheating of iron ore pellets","iron ore pellets, hot","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("preheating of hydrogen","hydrogen, hot",$null,$null)
    ,@("iron production, by electrowinning","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("leaching of iron ore","iron oxide in alkaline solution","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("market for cathode, graphite","cathode, graphite",$null,$null)
    ,@("nickel anode production, for electrolysis of iron ore","nickel anode, for electrolysis of iron ore","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("production of alkaline solution from sodium hydroxide of 50 wt-%","alkaline solution with 50 wt-% sodium hydroxide",$null,$null)
    ,@("steel production, electrowinning-electric arc furnace, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("steel production, electrowinning-electric arc furnace, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41122: Alloy steel in ingots or other primary forms and semi-finished products of alloy steel")
    ,@("ultrafine grinding of iron ore","iron ore, ultrafine ground","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
)

$startRow1 = 3109
for ($i = 0; $i -lt $newRowsPhase1.Count; $i++) {
    $r = $startRow1 + $i
    $rowVals = $newRowsPhase1[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $val = $rowVals[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

$filterLastRow = 3142

# Refresh the AutoFilter range (this also clears the now-stale sortState/sortCondition).
$ws.AutoFilterMode = $false
$ws.Range("A1:D" + $filterLastRow).AutoFilter()

# Keep the workbook-level _FilterDatabase name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$" + $filterLastRow
    }
}

$newRowsPhase2 = @(
    ,@("market for steel, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("market for steel, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, with carbon capture and storage, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, with carbon capture and storage, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, with top gas recycling, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, with top gas recycling, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, with top gas recycling, with carbon capture and storage, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, blast furnace-basic oxygen furnace, with top gas recycling, with carbon capture and storage, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("pig iron production, with natural gas-based direct reduction","iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("pig iron production, with natural gas-based direct reduction, with carbon capture and storage","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, natural gas-based direct reduction iron-electric arc furnace, with carbon capture and storage, low-alloyed","steel, low-alloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("steel production, natural gas-based direct reduction iron-electric arc furnace, with carbon capture and storage, unalloyed","steel, unalloyed","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
    ,@("pig iron production, by electrowinning","pig iron","2410:Manufacture of basic iron and steel","41111: Pig iron and spiegeleisen in pigs, blocks or other primary forms")
)

$startRow2 = 3143
for ($i = 0; $i -lt $newRowsPhase2.Count; $i++) {
    $r = $startRow2 + $i
    $rowVals = $newRowsPhase2[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $val = $rowVals[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

# Match the author's final selection / scroll position.
$ws.Range("A3138").Select()
$ws.Range("B3161").Select()
